$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2231155.5
$ws.Range("I116").Value = 8319
$ws.Range("K116").Value = 8319
$ws.Range("M116").Value = -4877
$ws.Range("H125").Value = 4000
$ws.Range("J125").Value = 7000
$ws.Range("L125").Value = 63000
$ws.Range("N125").Value = -67920
$ws.Range("H141").Value = 2418.0588
$ws.Range("I141").Value = 2418.0588
$ws.Range("K141").Value = 7254.176399999999
$ws.Range("M141").Value = -2074.176399999999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 43636.207
$ws.Range("I61").Value = 1812.1111
$ws.Range("J61").Value = 169108.5
$ws.Range("K61").Value = 1812.1111
$ws.Range("L61").Value = 169108.5
$ws.Range("M61").Value = -1600.1111
$ws.Range("N61").Value = -169532.5
$ws.Range("H97").Value = 544.8889
$ws.Range("I97").Value = 547.58826
$ws.Range("K97").Value = 547.58826
$ws.Range("M97").Value = -51.58825999999999
$ws.Range("H119").Value = 44333
$ws.Range("J119").Value = 44333
$ws.Range("L119").Value = 44333
$ws.Range("N119").Value = -54009
$ws.Range("H132").Value = 1735.6129
$ws.Range("I132").Value = 1735.6129
$ws.Range("K132").Value = 5206.8387
$ws.Range("M132").Value = -2676.8387
$ws.Range("H136").Value = 43636.207
$ws.Range("I136").Value = 1812.1111
$ws.Range("J136").Value = 169108.5
$ws.Range("K136").Value = 5436.3333
$ws.Range("L136").Value = 507325.5
$ws.Range("M136").Value = -2886.3333
$ws.Range("N136").Value = -512425.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5936
$ws.Range("I86").Value = 4581.3335
$ws.Range("K86").Value = 4581.3335
$ws.Range("M86").Value = -3458.3335
$ws.Range("H89").Value = 5936
$ws.Range("I89").Value = 4581.3335
$ws.Range("K89").Value = 22906.6675
$ws.Range("M89").Value = -17290.6675
$ws.Range("H94").Value = 1819.6666
$ws.Range("I94").Value = 1729.5
$ws.Range("K94").Value = 1729.5
$ws.Range("M94").Value = -1278.5
$ws.Range("H105").Value = 65282.375
$ws.Range("I105").Value = 101877.8
$ws.Range("K105").Value = 101877.8
$ws.Range("M105").Value = -100130.8
$ws.Range("H107").Value = 3202
$ws.Range("I107").Value = 2186.3333
$ws.Range("J107").Value = 4725.5
$ws.Range("K107").Value = 2186.3333
$ws.Range("L107").Value = 4725.5
$ws.Range("M107").Value = -266.3332999999998
$ws.Range("N107").Value = -8565.5
$ws.Range("H115").Value = 64106.168
$ws.Range("J115").Value = 64527.2
$ws.Range("L115").Value = 64527.2
$ws.Range("N115").Value = -67661.2
$ws.Range("H134").Value = 2534.125
$ws.Range("I134").Value = 1599.5
$ws.Range("K134").Value = 4798.5
$ws.Range("M134").Value = -2263.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 49999.5
$ws.Range("J50").Value = 49999.5
$ws.Range("L50").Value = 49999.5
$ws.Range("N50").Value = -51249.5
$ws.Range("H105").Value = 5641.6665
$ws.Range("I105").Value = 1283.3334
$ws.Range("J105").Value = 10000
$ws.Range("K105").Value = 1283.3334
$ws.Range("L105").Value = 10000
$ws.Range("M105").Value = 463.6666
$ws.Range("N105").Value = -13494
$ws.Range("H107").Value = 980
$ws.Range("I107").Value = 970.6667
$ws.Range("K107").Value = 970.6667
$ws.Range("M107").Value = 949.3333
$ws.Range("H132").Value = 1346.6086
$ws.Range("I132").Value = 1346.6086
$ws.Range("K132").Value = 4039.8258
$ws.Range("M132").Value = -1509.8258
$ws.Range("H134").Value = 45969.87
$ws.Range("I134").Value = 2464.1052
$ws.Range("K134").Value = 7392.3156
$ws.Range("M134").Value = -4857.3156

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 106.22727
$ws.Range("J2").Value = 136.375
$ws.Range("L2").Value = 818.25
$ws.Range("N2").Value = -1044.25
$ws.Range("H4").Value = 71798020
$ws.Range("I4").Value = 78908820
$ws.Range("K4").Value = 236726460
$ws.Range("M4").Value = -236726348
$ws.Range("H7").Value = 6736.1665
$ws.Range("J7").Value = 20010
$ws.Range("L7").Value = 60030
$ws.Range("N7").Value = -60254
$ws.Range("H17").Value = 810
$ws.Range("I17").Value = 20
$ws.Range("J17").Value = 941.6667
$ws.Range("K17").Value = 60
$ws.Range("L17").Value = 2825.0001
$ws.Range("M17").Value = 109
$ws.Range("N17").Value = -3163.0001
$ws.Range("H122").Value = 2020868.4
$ws.Range("I122").Value = 691.5
$ws.Range("J122").Value = 3367653
$ws.Range("K122").Value = 6223.5
$ws.Range("L122").Value = 30308877
$ws.Range("M122").Value = -3773.5
$ws.Range("N122").Value = -30313777
$ws.Range("H131").Value = 41404.32
$ws.Range("I131").Value = 111648.664
$ws.Range("K131").Value = 334945.992
$ws.Range("M131").Value = -329905.992
$ws.Range("H132").Value = 2554.6191
$ws.Range("I132").Value = 4182.6665
$ws.Range("J132").Value = 2283.2778
$ws.Range("K132").Value = 37643.9985
$ws.Range("L132").Value = 20549.5002
$ws.Range("M132").Value = -35113.9985
$ws.Range("N132").Value = -25609.5002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 31354.8
$ws.Range("J44").Value = 31354.8
$ws.Range("L44").Value = 31354.8
$ws.Range("N44").Value = -32546.8
$ws.Range("H47").Value = 29963.334
$ws.Range("J47").Value = 29945
$ws.Range("L47").Value = 29945
$ws.Range("N47").Value = -31081
$ws.Range("H97").Value = 371156.56
$ws.Range("I97").Value = 500321.2
$ws.Range("J97").Value = 2114.7144
$ws.Range("K97").Value = 500321.2
$ws.Range("L97").Value = 2114.7144
$ws.Range("M97").Value = -499825.2
$ws.Range("N97").Value = -3106.7144
$ws.Range("H113").Value = 3032344.8
$ws.Range("I113").Value = 1039.8
$ws.Range("J113").Value = 5558432.5
$ws.Range("K113").Value = 1039.8
$ws.Range("L113").Value = 5558432.5
$ws.Range("M113").Value = 1130.2
$ws.Range("N113").Value = -5562772.5
$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -53494
$ws.Range("H132").Value = 5327.037
$ws.Range("I132").Value = 4003.6365
$ws.Range("J132").Value = 11150
$ws.Range("K132").Value = 12010.9095
$ws.Range("L132").Value = 33450
$ws.Range("M132").Value = -9480.9095
$ws.Range("N132").Value = -38510

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 70497.42999999999
$ws.Range("I7").Value = 37247
$ws.Range("K7").Value = 37247
$ws.Range("M7").Value = -37135
$ws.Range("H55").Value = 1671.6444
$ws.Range("I55").Value = 932.2
$ws.Range("J55").Value = 3150.5334
$ws.Range("K55").Value = 932.2
$ws.Range("L55").Value = 3150.5334
$ws.Range("M55").Value = -759.2
$ws.Range("N55").Value = -3496.5334
$ws.Range("H61").Value = 1073.5
$ws.Range("I61").Value = 1073.5
$ws.Range("K61").Value = 1073.5
$ws.Range("M61").Value = -871.5
$ws.Range("H93").Value = 1392.2
$ws.Range("I93").Value = 1104.1666
$ws.Range("K93").Value = 1104.1666
$ws.Range("M93").Value = 143.8334
$ws.Range("H113").Value = 1073.5
$ws.Range("I113").Value = 1073.5
$ws.Range("K113").Value = 1073.5
$ws.Range("M113").Value = 1096.5
$ws.Range("H126").Value = 70497.42999999999
$ws.Range("I126").Value = 37247
$ws.Range("K126").Value = 111741
$ws.Range("M126").Value = -109271
$ws.Range("H136").Value = 6424.96
$ws.Range("I136").Value = 6590.5835
$ws.Range("J136").Value = 6272.077
$ws.Range("K136").Value = 19771.7505
$ws.Range("L136").Value = 18816.231
$ws.Range("M136").Value = -17221.7505
$ws.Range("N136").Value = -23916.231

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 7125
$ws.Range("I55").Value = 4000
$ws.Range("J55").Value = 8166.6665
$ws.Range("K55").Value = 4000
$ws.Range("L55").Value = 8166.6665
$ws.Range("M55").Value = -3723
$ws.Range("N55").Value = -8720.666499999999
$ws.Range("H61").Value = 1465664.6
$ws.Range("I61").Value = 2535038
$ws.Range("J61").Value = 39833.332
$ws.Range("K61").Value = 2535038
$ws.Range("L61").Value = 39833.332
$ws.Range("M61").Value = -2534746
$ws.Range("N61").Value = -40417.332
$ws.Range("H81").Value = 33712.715
$ws.Range("I81").Value = 2568.5715
$ws.Range("K81").Value = 5137.143
$ws.Range("M81").Value = -4076.143
$ws.Range("H84").Value = 33712.715
$ws.Range("I84").Value = 2568.5715
$ws.Range("K84").Value = 25685.715
$ws.Range("M84").Value = -20381.715
$ws.Range("H100").Value = 5954215.5
$ws.Range("I100").Value = 10206538
$ws.Range("J100").Value = 964
$ws.Range("K100").Value = 20413076
$ws.Range("L100").Value = 1928
$ws.Range("M100").Value = -20412535
$ws.Range("N100").Value = -3010
$ws.Range("H113").Value = 1592.7142
$ws.Range("I113").Value = 1525
$ws.Range("K113").Value = 4575
$ws.Range("M113").Value = -2405
